$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Title paragraphs: "Introduction to using R for Spatial Analysis"
#    -> "Introduction to using " + <line break> + "R as a GIS"
#    (same character formatting, just split into two runs around a
#    manual line break).
# ------------------------------------------------------------------
$titleXmlTemplate = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
{0}
<w:r>
<w:rPr>
<w:rFonts w:asciiTheme="majorHAnsi" w:eastAsia="Times New Roman" w:hAnsiTheme="majorHAnsi" w:cs="Times New Roman"/>
<w:b/>
<w:bCs/>
<w:color w:val="222222"/>
<w:sz w:val="72"/>
<w:szCs w:val="36"/>
<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
</w:rPr>
<w:t xml:space="preserve">Introduction to using </w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:asciiTheme="majorHAnsi" w:eastAsia="Times New Roman" w:hAnsiTheme="majorHAnsi" w:cs="Times New Roman"/>
<w:b/>
<w:bCs/>
<w:color w:val="222222"/>
<w:sz w:val="72"/>
<w:szCs w:val="36"/>
<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
</w:rPr>
<w:br/>
<w:t>R as a GIS</w:t>
</w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$oldTitle = "Introduction to using R for Spatial Analysis"

$targets = New-Object System.Collections.ArrayList
foreach ($p in $d.Paragraphs) {
    $ptext = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($ptext -eq $oldTitle) {
        [void]$targets.Add($p)
    }
}

foreach ($p in $targets) {
    $pPrXml = ""
    $wopx = $p.Range.WordOpenXML
    if ($wopx -match "(?s)<w:p[ >].*?(<w:pPr>.*?</w:pPr>)") {
        $pPrXml = $matches[1]
    }
    $xml = $titleXmlTemplate -f $pPrXml
    $p.Range.InsertXML($xml)
}

# ------------------------------------------------------------------
# 2) Date paragraph: "Tuesday 1st December 2015" (with a couple of
#    whitespace variants across the document) -> "Thursday 3rd March
#    2016", collapsed into a single run.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Tuesday 1st December  2015", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Thursday 3rd March 2016", 2) | Out-Null
$d.Content.Find.Execute("Tuesday 1st December 2015", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Thursday 3rd March 2016", 2) | Out-Null
